# Apply "Nogle flere Chancekort funktioner er implementeret" edits.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ChanceKort")

# Reword four chance-card texts (C11, C3, C12, C19) - order matters for
# the order in which new shared-string entries get appended.
$ws.Range("C11").Value = "Du løslades uden omkostninger. Du har nu dette indtil du får brug for det."
$ws.Range("C3").Value = "Du rykkes frem til start."
$ws.Range("C12").Value = "Du rykkes frem til Strandpromenaden."
$ws.Range("C19").Value = "Gratis felt. Du kykkes frem til Skaterparken for at lave det perfekte grind. Hvis ingen ejer den, får du den gratis. Ellers skal du betale leje ejeren."

# Update the saved view state (scroll position & active selection).
$ws.Activate()
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("C19").Select()
